# Edit script: add "2022-Q3" sheet with fund holding data, and insert a new
# leading row into the "总计" (summary) sheet for the 2022-Q3 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" summary sheet: a new 2022-Q3 row is inserted above the
#    existing rows and every subsequent row shifts down by one.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Carry the bold/bordered "index" style (column A) down onto the new row 8.
$zj.Cells.Item(7, 1).Copy()
$zj.Range("A8").PasteSpecial(-4122)

$summaryData = @(
    @(0, '2022-Q3', 29, 9.470000000000001),
    @(1, '2022-Q2', 12, 0.51),
    @(2, '2022-Q1', 2, 0.07000000000000001),
    @(3, '2021-Q4', 8, 0.4),
    @(4, '2021-Q3', 4, 0.21),
    @(5, '2021-Q2', 2, 0.01),
    @(6, '2020-Q4', 2, 0.08)
)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $r = $i + 2
    $row = $summaryData[$i]
    $zj.Cells.Item($r, 1).Value = $row[0]
    $zj.Cells.Item($r, 2).Value = $row[1]
    $zj.Cells.Item($r, 3).Value = $row[2]
    $zj.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" worksheet right after "总计" and populate it
#    with the quarterly fund-holding table.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $zj)
$q3.Name = "2022-Q3"

# Header row (B1:H1) — bold/centered/bordered, matching the other quarter sheets.
$zj.Cells.Item(1, 2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$headers = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Column A (row index 0..28) — bold/centered/bordered, same style as "总计"!A2.
$zj.Cells.Item(2, 1).Copy()
$q3.Range("A2:A30").PasteSpecial(-4122)

# Columns D:G hold numeric-looking text (e.g. "75.87") in the source data, so
# force them to Text before writing — otherwise they'd be auto-coerced into
# numbers and lose things like trailing zeros.
$q3.Range("D2:G30").NumberFormat = "@"

$q3Data = @(
    @(0, '007130', '中庚小盘价值股票', '75.87', '93.06', '5.83', '4.4232', 1),
    @(1, '007497', '中庚价值灵动灵活配置混合', '36.46', '89.30', '3.91', '1.4256', 4),
    @(2, '015182', '汇添富逆向投资混合D', '21.25', '92.44', '4.01', '0.8521', 9),
    @(3, '470098', '汇添富逆向投资混合A', '20.96', '92.44', '4.01', '0.8405', 9),
    @(4, '013552', '汇添富成长领先混合A', '10.27', '62.59', '3.25', '0.3338', 8),
    @(5, '009989', '华宝研究精选混合', '8.14', '82.83', '3.01', '0.2450', 10),
    @(6, '005416', '鹏华尊惠18个月定期开放混合A', '7.53', '37.52', '3.16', '0.2379', 1),
    @(7, '010296', '万家互联互通中国优势量化策略混合A', '4.22', '86.53', '5.35', '0.2258', 4),
    @(8, '013553', '汇添富成长领先混合C', '5.87', '62.59', '3.25', '0.1908', 8),
    @(9, '020015', '国泰区位优势混合A', '1.95', '87.12', '5.64', '0.1100', 3),
    @(10, '000866', '华宝高端制造股票', '2.30', '89.16', '3.18', '0.0731', 10),
    @(11, '009667', '鹏华安庆混合A', '2.12', '38.12', '3.05', '0.0647', 1),
    @(12, '003165', '鹏华弘嘉灵活配置混合A', '0.82', '91.09', '6.99', '0.0573', 1),
    @(13, '009230', '鹏华安和混合A', '1.80', '39.30', '3.04', '0.0547', 1),
    @(14, '011572', '鹏华安荣混合A', '1.50', '39.92', '3.13', '0.0470', 1),
    @(15, '014509', '汇添富先进制造混合C', '0.85', '87.50', '4.88', '0.0415', 6),
    @(16, '014508', '汇添富先进制造混合A', '0.83', '87.50', '4.88', '0.0405', 6),
    @(17, '009668', '鹏华安庆混合C', '1.32', '38.12', '3.05', '0.0403', 1),
    @(18, '010297', '万家互联互通中国优势量化策略混合C', '0.46', '86.53', '5.35', '0.0246', 4),
    @(19, '005482', '博时创新驱动灵活配置混合A', '0.39', '85.84', '5.98', '0.0233', 2),
    @(20, '014141', '大成新能源混合A', '0.52', '82.56', '4.46', '0.0232', 9),
    @(21, '009231', '鹏华安和混合C', '0.73', '39.30', '3.04', '0.0222', 1),
    @(22, '005417', '鹏华尊惠18个月定期开放混合C', '0.58', '37.52', '3.16', '0.0183', 1),
    @(23, '003166', '鹏华弘嘉灵活配置混合C', '0.25', '91.09', '6.99', '0.0175', 1),
    @(24, '014142', '大成新能源混合C', '0.39', '82.56', '4.46', '0.0174', 9),
    @(25, '011573', '鹏华安荣混合C', '0.34', '39.92', '3.13', '0.0106', 1),
    @(26, '015181', '汇添富逆向投资混合C', '0.24', '92.44', '4.01', '0.0096', 9),
    @(27, '005483', '博时创新驱动灵活配置混合C', '0.06', '85.84', '5.98', '0.0036', 2),
    @(28, '015594', '国泰区位优势混合C', '0.00', '87.12', '5.64', $null, 3)
)


for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    if ($row[6] -eq $null) {
        $q3.Cells.Item($r, 7).NumberFormat = "General"
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        $q3.Cells.Item($r, 7).Value = $row[6]
    }
    $q3.Cells.Item($r, 8).Value = $row[7]
}
